# "Updated translation patch to support the huge v2 update."
#
# This introduces a dedicated translation column (B) for every row of the
# sheet: for rows that did not yet have a translation, column B is seeded
# with a copy of the original (column A) text; for the handful of rows
# that already had translated text stashed away in columns C/D, that text
# is moved into column B instead. Columns C and D (now unused) are then
# removed, shrinking the sheet's dimension from A1:D47 down to A1:B47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-41, 45 and 46: column B does not yet have a translation, so it is
# simply initialized with a copy of column A's text for that row.
$plainCopyRows = (1..41) + (45, 46)
foreach ($r in $plainCopyRows) {
    $ws.Range("A$r").Copy()
    $ws.Range("B$r").PasteSpecial()
}

# Rows 42-44 already had translated text living in column D (and, for 43
# and 44, duplicated into column C as well). Bring that translation into
# column B.
foreach ($r in 42, 43, 44) {
    $ws.Range("D$r").Copy()
    $ws.Range("B$r").PasteSpecial()
}

# Row 47 already carries its translation in column B ("cake"), so it is
# left untouched.

# Finally, drop the now-redundant columns C and D so the sheet is left
# with only A (original) and B (translation).
$ws.Columns.Item(3).Delete()
$ws.Columns.Item(3).Delete()
